$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Neurology" to "Session"
$ws.Name = "Session"

# New QR-scanner backup rows to append below the existing log (rows 76-85).
# Every column in this sheet is stored as text (numeric-looking IDs/dates/
# times included), so force text formatting before writing the values to
# keep them as strings rather than being coerced into numbers/dates.
$newRows = @(
    @(76, "201888", "Neurology", "29/12/2025", "11:00:55", "Scan", "emp17.farah.a.youssef@gmail.com"),
    @(77, "201243", "Neurology", "29/12/2025", "11:06:16", "Scan", "emp17.farah.a.youssef@gmail.com"),
    @(78, "201479", "Neurology", "29/12/2025", "11:10:24", "Scan", "emp17.farah.a.youssef@gmail.com"),
    @(79, "201322", "Neurology", "29/12/2025", "13:39:28", "Scan", "emp17.farah.a.youssef@gmail.com"),
    @(80, "201265", "Neurology", "29/12/2025", "13:39:32", "Scan", "emp17.farah.a.youssef@gmail.com"),
    @(81, "201234", "Neurology", "29/12/2025", "13:39:37", "Scan", "emp17.farah.a.youssef@gmail.com"),
    @(82, "201329", "Neurology", "29/12/2025", "13:39:39", "Scan", "emp17.farah.a.youssef@gmail.com"),
    @(83, "201239", "Neurology", "29/12/2025", "13:39:45", "Scan", "emp17.farah.a.youssef@gmail.com"),
    @(84, "201061", "Neurology", "29/12/2025", "13:40:04", "Scan", "emp17.farah.a.youssef@gmail.com"),
    @(85, "201416", "Neurology", "29/12/2025", "13:40:09", "Scan", "emp17.farah.a.youssef@gmail.com")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $rowRange = $ws.Range("A" + $r + ":F" + $r)
    $rowRange.NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    # Values are now stored as text; restore the default "Normal" style so
    # the new cells match the rest of the sheet (which carries no explicit
    # per-cell formatting) instead of keeping the temporary "@" text format.
    $rowRange.Style = "Normal"
}
